# river update May 2024
# Applies the recomputed statistic values for the existing RepSite/data rows
# (recalculated Mean/83rd/92nd percentile figures etc. after new raw data was
# folded into the underlying calculations) and appends three brand-new rows
# (316-318: ASPM / MCI / QMCI for the "2019 - 2023" period) that extend the
# sheets used range from U315 to U318.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated computed values ---
$ws.Range("G8").Value = 1897.49112569783
$ws.Range("G9").Value = 1897.49112569783
$ws.Range("G10").Value = 1897.49112569783
$ws.Range("G11").Value = 1897.49112569783
$ws.Range("G18").Value = 1734.25472080301
$ws.Range("G19").Value = 1734.25472080301
$ws.Range("G20").Value = 1734.25472080301
$ws.Range("G21").Value = 1734.25472080301
$ws.Range("G27").Value = 0.239645007123736
$ws.Range("L27").Value = 0.08624999999999999
$ws.Range("G28").Value = 0.239645007123736
$ws.Range("L28").Value = 0.08624999999999999
$ws.Range("G32").Value = 1643.43650526397
$ws.Range("G33").Value = 1643.43650526397
$ws.Range("G34").Value = 1643.43650526397
$ws.Range("G35").Value = 1643.43650526397
$ws.Range("G41").Value = 0.241744886549602
$ws.Range("L41").Value = 0.08624999999999999
$ws.Range("G42").Value = 0.241744886549602
$ws.Range("L42").Value = 0.08624999999999999
$ws.Range("G46").Value = 2434.33983832991
$ws.Range("H46").Value = 36942.9733194293
$ws.Range("M46").Value = 3841.93741
$ws.Range("G47").Value = 2434.33983832991
$ws.Range("H47").Value = 36942.9733194293
$ws.Range("M47").Value = 3841.93741
$ws.Range("G48").Value = 2434.33983832991
$ws.Range("H48").Value = 36942.9733194293
$ws.Range("M48").Value = 3841.93741
$ws.Range("G49").Value = 2434.33983832991
$ws.Range("H49").Value = 36942.9733194293
$ws.Range("M49").Value = 3841.93741
$ws.Range("F55").Value = 0.2425
$ws.Range("G55").Value = 0.252358534021241
$ws.Range("F56").Value = 0.2425
$ws.Range("G56").Value = 0.252358534021241
$ws.Range("G60").Value = 3361.62821411273
$ws.Range("I60").Value = 21665
$ws.Range("G61").Value = 3361.62821411273
$ws.Range("I61").Value = 21665
$ws.Range("G62").Value = 3361.62821411273
$ws.Range("I62").Value = 21665
$ws.Range("G63").Value = 3361.62821411273
$ws.Range("I63").Value = 21665
$ws.Range("F69").Value = 0.24925
$ws.Range("G69").Value = 0.261836559785098
$ws.Range("F70").Value = 0.24925
$ws.Range("G70").Value = 0.261836559785098
$ws.Range("G74").Value = 2631.44173359954
$ws.Range("G75").Value = 2631.44173359954
$ws.Range("G76").Value = 2631.44173359954
$ws.Range("G77").Value = 2631.44173359954
$ws.Range("F83").Value = 0.24075
$ws.Range("G83").Value = 0.265954765234348
$ws.Range("F84").Value = 0.24075
$ws.Range("G84").Value = 0.265954765234348
$ws.Range("G88").Value = 2779.31337207644
$ws.Range("G89").Value = 2779.31337207644
$ws.Range("G90").Value = 2779.31337207644
$ws.Range("G91").Value = 2779.31337207644
$ws.Range("F97").Value = 0.2375
$ws.Range("G97").Value = 0.263079555346172
$ws.Range("F98").Value = 0.2375
$ws.Range("G98").Value = 0.263079555346172
$ws.Range("G102").Value = 2669.31337207644
$ws.Range("G103").Value = 2669.31337207644
$ws.Range("G104").Value = 2669.31337207644
$ws.Range("G105").Value = 2669.31337207644
$ws.Range("G111").Value = 0.255308794666438
$ws.Range("G112").Value = 0.255308794666438
$ws.Range("G129").Value = 0.254281516130486
$ws.Range("M129").Value = 0.44789
$ws.Range("G130").Value = 0.254281516130486
$ws.Range("M130").Value = 0.44789
$ws.Range("G147").Value = 0.267367124276594
$ws.Range("H147").Value = 1.1451
$ws.Range("G148").Value = 0.267367124276594
$ws.Range("H148").Value = 1.1451
$ws.Range("G167").Value = 0.253618
$ws.Range("H167").Value = 1.1451
$ws.Range("M167").Value = 0.4928
$ws.Range("G168").Value = 0.253618
$ws.Range("H168").Value = 1.1451
$ws.Range("M168").Value = 0.4928
$ws.Range("G174").Value = 0.899330792360192
$ws.Range("H174").Value = 4.79555486384846
$ws.Range("G187").Value = 0.251157423097401
$ws.Range("H187").Value = 1.1451
$ws.Range("M187").Value = 0.5000599999999999
$ws.Range("N187").Value = 0.58678
$ws.Range("G188").Value = 0.251157423097401
$ws.Range("H188").Value = 1.1451
$ws.Range("M188").Value = 0.5000599999999999
$ws.Range("N188").Value = 0.58678
$ws.Range("G194").Value = 0.865901219196632
$ws.Range("H194").Value = 4.79555486384846
$ws.Range("G197").Value = 1037.51622562838
$ws.Range("H197").Value = 12458.4250164957
$ws.Range("G198").Value = 1037.51622562838
$ws.Range("H198").Value = 12458.4250164957
$ws.Range("G199").Value = 1037.51622562838
$ws.Range("H199").Value = 12458.4250164957
$ws.Range("G200").Value = 1037.51622562838
$ws.Range("H200").Value = 12458.4250164957
$ws.Range("G207").Value = 0.276043396743745
$ws.Range("H207").Value = 1.1451
$ws.Range("M207").Value = 0.50926
$ws.Range("G208").Value = 0.276043396743745
$ws.Range("H208").Value = 1.1451
$ws.Range("M208").Value = 0.50926
$ws.Range("G214").Value = 0.919911097276969
$ws.Range("H214").Value = 4.79555486384846
$ws.Range("G217").Value = 893.839933051449
$ws.Range("H217").Value = 12458.4250164957
$ws.Range("G218").Value = 893.839933051449
$ws.Range("H218").Value = 12458.4250164957
$ws.Range("G219").Value = 893.839933051449
$ws.Range("H219").Value = 12458.4250164957
$ws.Range("G220").Value = 893.839933051449
$ws.Range("H220").Value = 12458.4250164957
$ws.Range("G227").Value = 0.28275514779232
$ws.Range("H227").Value = 1.1451
$ws.Range("M227").Value = 0.51577
$ws.Range("N227").Value = 0.6086
$ws.Range("G228").Value = 0.28275514779232
$ws.Range("H228").Value = 1.1451
$ws.Range("M228").Value = 0.51577
$ws.Range("N228").Value = 0.6086
$ws.Range("G234").Value = 0.927777763060605
$ws.Range("H234").Value = 4.79555486384846
$ws.Range("G237").Value = 718.186343093068
$ws.Range("H237").Value = 12458.4250164957
$ws.Range("G238").Value = 718.186343093068
$ws.Range("H238").Value = 12458.4250164957
$ws.Range("G239").Value = 718.186343093068
$ws.Range("H239").Value = 12458.4250164957
$ws.Range("G240").Value = 718.186343093068
$ws.Range("H240").Value = 12458.4250164957
$ws.Range("F247").Value = 0.2591
$ws.Range("G247").Value = 0.29528383819174
$ws.Range("L247").Value = 0.07065
$ws.Range("F248").Value = 0.2591
$ws.Range("G248").Value = 0.29528383819174
$ws.Range("L248").Value = 0.07065
$ws.Range("G254").Value = 0.984558065087465
$ws.Range("H254").Value = 4.79555486384846
$ws.Range("G257").Value = 767.997365615587
$ws.Range("H257").Value = 12458.4250164957
$ws.Range("G258").Value = 767.997365615587
$ws.Range("H258").Value = 12458.4250164957
$ws.Range("G259").Value = 767.997365615587
$ws.Range("H259").Value = 12458.4250164957
$ws.Range("G260").Value = 767.997365615587
$ws.Range("H260").Value = 12458.4250164957
$ws.Range("F267").Value = 0.2353
$ws.Range("G267").Value = 0.293020072854515
$ws.Range("L267").Value = 0.07065
$ws.Range("F268").Value = 0.2353
$ws.Range("G268").Value = 0.293020072854515
$ws.Range("L268").Value = 0.07065
$ws.Range("G277").Value = 887.273961360268
$ws.Range("H277").Value = 12458.4250164957
$ws.Range("G278").Value = 887.273961360268
$ws.Range("H278").Value = 12458.4250164957
$ws.Range("G279").Value = 887.273961360268
$ws.Range("H279").Value = 12458.4250164957
$ws.Range("G280").Value = 887.273961360268
$ws.Range("H280").Value = 12458.4250164957
$ws.Range("G287").Value = 0.306489361702128
$ws.Range("G288").Value = 0.306489361702128
$ws.Range("F307").Value = 0.3505
$ws.Range("G307").Value = 0.309179487179487
$ws.Range("F308").Value = 0.3505
$ws.Range("G308").Value = 0.309179487179487

# --- New rows 316-318 ---
# Row 316
$ws.Range("A316").Value = "Whanganui at Wades Landing"
$ws.Range("B316").Value = "ASPM"
$ws.Range("C316").Value = "C"
$ws.Range("D316").Value = "2019 - 2023"
$ws.Range("E316").Value = "RepSite"
$ws.Range("F316").Value = 0.362
$ws.Range("G316").Value = 0.4006
$ws.Range("H316").Value = 0.545
$ws.Range("I316").Value = 0.545
$ws.Range("L316").Value = 0.344
$ws.Range("M316").Value = 0.5257500000000001
$ws.Range("N316").Value = 0.545
$ws.Range("O316").Value = 1778217
$ws.Range("P316").Value = 5668812
$ws.Range("Q316").Value = "Ruapehu District"
$ws.Range("R316").Value = "Whanganui"
$ws.Range("S316").Value = "Pipiriki"
$ws.Range("T316").Value = "Whai_5a"

# Row 317
$ws.Range("A317").Value = "Whanganui at Wades Landing"
$ws.Range("B317").Value = "MCI"
$ws.Range("C317").Value = "C"
$ws.Range("D317").Value = "2019 - 2023"
$ws.Range("E317").Value = "RepSite"
$ws.Range("F317").Value = 103
$ws.Range("G317").Value = 97.30200000000001
$ws.Range("H317").Value = 108.18
$ws.Range("I317").Value = 108.18
$ws.Range("L317").Value = 100
$ws.Range("M317").Value = 107.067
$ws.Range("N317").Value = 108.18
$ws.Range("O317").Value = 1778217
$ws.Range("P317").Value = 5668812
$ws.Range("Q317").Value = "Ruapehu District"
$ws.Range("R317").Value = "Whanganui"
$ws.Range("S317").Value = "Pipiriki"
$ws.Range("T317").Value = "Whai_5a"

# Row 318
$ws.Range("A318").Value = "Whanganui at Wades Landing"
$ws.Range("B318").Value = "QMCI"
$ws.Range("C318").Value = "D"
$ws.Range("D318").Value = "2019 - 2023"
$ws.Range("E318").Value = "RepSite"
$ws.Range("F318").Value = 4.133
$ws.Range("G318").Value = 4.4612
$ws.Range("H318").Value = 5.57
$ws.Range("I318").Value = 5.57
$ws.Range("L318").Value = 3.9265
$ws.Range("M318").Value = 5.55005
$ws.Range("N318").Value = 5.57
$ws.Range("O318").Value = 1778217
$ws.Range("P318").Value = 5668812
$ws.Range("Q318").Value = "Ruapehu District"
$ws.Range("R318").Value = "Whanganui"
$ws.Range("S318").Value = "Pipiriki"
$ws.Range("T318").Value = "Whai_5a"

